# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns, and reorders a few rows
# (B/C/D/E) for Aptos/PEPE and PancakeSwap/Monero swaps, per commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text is unambiguous (keeps inline/shared string type automatically) ---
$ws.Range("D2").Value = "58.815.05"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.520.20"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").Value = "2.519.12"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "2.963.82"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "58.800.15"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "2.511.95"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("E24").Value = "  +5.51%  "
$ws.Range("E25").Value = "  +4.72%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0769"
$ws.Range("E29").Value = "  +3.17%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +10.96%  "
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  +4.54%  "
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("E45").Value = "  +4.80%  "
$ws.Range("E46").Value = "  +9.20%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +3.34%  "

# --- Cells whose new text would be auto-parsed as a number by Excel; force them to stay text ---
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "535.81"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "136.37"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.566"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.349"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "23.05"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.08"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.27"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "322.74"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.98"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.20"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.421"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.52"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.65"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "170.59"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.75"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "18.35"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "284.17"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.06"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.608"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "130.16"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.86"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0923"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0219"
$c.Style = "Normal"
